# Update "想去人数" (F column) counts on sheets 展览, 演出, and 全部类型
# to reflect the newer scrape of 广州-漫展信息.

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 121
$ws1.Range("F4").Value  = 438
$ws1.Range("F6").Value  = 136
$ws1.Range("F7").Value  = 1214
$ws1.Range("F8").Value  = 407
$ws1.Range("F12").Value = 384
$ws1.Range("F16").Value = 733
$ws1.Range("F17").Value = 296
$ws1.Range("F19").Value = 1027
$ws1.Range("F20").Value = 477
$ws1.Range("F23").Value = 390
$ws1.Range("F25").Value = 48
$ws1.Range("F26").Value = 478

# Sheet: 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 287

# Sheet: 全部类型 (All types combined)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 121
$ws4.Range("F6").Value  = 438
$ws4.Range("F8").Value  = 136
$ws4.Range("F9").Value  = 1214
$ws4.Range("F10").Value = 407
$ws4.Range("F17").Value = 384
$ws4.Range("F19").Value = 287
$ws4.Range("F23").Value = 733
$ws4.Range("F24").Value = 296
$ws4.Range("F26").Value = 1027
$ws4.Range("F27").Value = 477
$ws4.Range("F32").Value = 390
$ws4.Range("F36").Value = 48
$ws4.Range("F38").Value = 478
